$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.716.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.319.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.60%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '97.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '272.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.656.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.53%  '
$ws.Range('E16').Value = '  +8.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.325.30'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.727.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  +4.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.58%  '
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('E24').Value = '  +3.66%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.29'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.44'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '175.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0916'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E35').Value = '  +2.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0361'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('E38').Value = '  +2.42%  '
$ws.Range('E39').Value = '  -5.75%  '
$ws.Range('E40').Value = '  +5.04%  '
$ws.Range('E41').Value = '  +8.08%  '
$ws.Range('E42').Value = '  +20.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.80%  '
$ws.Range('E44').Value = '  +10.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +16.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.542.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.44%  '
